$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(5.463630573248407, 2.240568743774976, 1.322903771964442, 0.7159837955602021, 0.6407549044741776, -0.4845913390137646, 25.35414012738854)
    3 = @(4.633661436405054, 1.683809642222656, 0.8144442427531886, 0.5966637052870977, 0.5497627707135105, 0.8797663191006361, 25.17899060047703)
    4 = @(5.231111111111111, 2.311781089734088, 1.295143358830755, 0.7691535895238956, 0.5403056057900927, -0.4549955092976475, 24.35636856368564)
    5 = @(4.720194510067887, 1.7711033548038, 0.8911236331312969, 0.6063470312657683, 0.4324091930454509, 0.8113933948058411, 24.21593482427432)
    6 = @(5.521489932885906, 2.443767023738328, 1.350303711758079, 0.7956553141680518, 0.6276022584510498, -0.5169876928642829, 25.55704697986577)
    7 = @(4.749391535190125, 1.789747079028432, 0.829383928103861, 0.5885783525344821, 0.5083603725834306, 0.815415571167421, 25.18164722291512)
    8 = @(5.100556492411467, 2.258784850519582, 1.275894941232602, 0.7613507843784619, 0.5336702779734646, -0.4173430637859449, 29.81281618887015)
    9 = @(4.484853849112721, 1.697809794214145, 0.8562049726558092, 0.5345199867736541, 0.4178948825234093, 0.6975964498132601, 26.98104385230003)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    $ws.Range("F$row").Value = $values[4]
    $ws.Range("G$row").Value = $values[5]
    $ws.Range("H$row").Value = $values[6]
}
